$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14: categorie accessibilite/SEO + new "erreur 404" problem
$ws.Range("A14").Value = "accessibilité/SEO"
$ws.Range("B14").Value = "erreur 404"

# Row 15: categorie accessibilite + new "aria" problem
$ws.Range("A15").Value = "accessibilité"
$ws.Range("B15").Value = "aria"

# Update selected cell to A15 (as stored in the saved view state)
$ws.Range("A15").Select()
